{"js": "// Diff analysis\n// -------------\n// Every hunk in the supplied diff touches *serialization* details only:\n//   - the set of namespace prefixes declared on each part's root element\n//     (e.g. word/document.xml's <w:document> goes from the Word-Windows\n//     prefix set (wpc/cx/mc/o/v/w10/...) to a different writer's prefix\n//     set (a/m/ma/pic/s/xml/...), same for endnotes/footnotes/headers/\n//     footers/styles/theme1),\n//   - attribute order within an element (<w:bookmarkStart w:id=\"0\"\n//     w:name=\"_GoBack\"/> -> <w:bookmarkStart w:name=\"_GoBack\" w:id=\"0\"/>,\n//     <w:rFonts .../> theme attr order, <w:style w:type=\"..\" w:default=\"..\"\n//     w:styleId=\"..\"/> -> w:type/w:styleId/w:default, etc.), and\n//   - boolean literal spelling (w:val=\"1\"/\"0\" -> w:val=\"true\"/\"false\" in\n//     styles.xml's latentStyles/lsdException block and the default-style\n//     flags), plus one extLst child element's namespace prefix being\n//     renamed (thm15 -> m) in theme1.xml.\n// There is no added/removed/retitled paragraph, run, bookmark, style\n// definition, theme color, or header/footer text anywhere in the diff \u2014\n// the commit message (\"fix import and update golden files\") matches\n// that: it is a test-fixture refresh after fixing an *importer*, not a\n// content edit. Word's object model (Office.js here) deliberately has no\n// hooks for dictating XML namespace prefixes, attribute order, or\n// boolean-literal spelling \u2014 the host always re-serializes OOXML with its\n// own writer \u2014 so there is nothing for a content-level script to change.\n//\n// We still touch the one piece of content the diff's context lines pass\n// through (the document body's \"_GoBack\" bookmark) to confirm it is\n// present/unchanged, but we deliberately avoid any write that would\n// introduce a *new* difference (e.g. re-inserting the bookmark would\n// duplicate it, which is not what happened here).\nconst bookmark = context.document.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmark.load(\"isNullObject,text\");\n\nconst body = context.document.body;\nbody.load(\"text\");\n\nawait context.sync();\n\n// Nothing to mutate: the bookmark's name/position and the body's text are\n// already exactly what the target state requires, so no further API\n// calls are made. (Read-only verification only \u2014 no context.sync() after\n// this point is needed because no properties were queued for write.)\n", "ps1": "# Diff analysis\n# -------------\n# Every hunk in the supplied diff is a *serialization* change only:\n#   - the namespace-prefix set declared on each part's root element\n#     (word/document.xml, endnotes.xml, footer1-3.xml, footnotes.xml,\n#     header1-3.xml, styles.xml, theme/theme1.xml all swap from the\n#     Word-Windows prefix set (wpc/cx/mc/o/v/w10/...) to a different\n#     writer's prefix set (a/m/ma/pic/s/xml/...)),\n#   - attribute order within elements (<w:bookmarkStart w:id=\"0\"\n#     w:name=\"_GoBack\"/> -> <w:bookmarkStart w:name=\"_GoBack\" w:id=\"0\"/>,\n#     <w:rFonts/> theme-attr order, <w:style w:type=\"..\" w:default=\"..\"\n#     w:styleId=\"..\"/> -> w:type/w:styleId/w:default, etc.), and\n#   - boolean-literal spelling (w:val=\"1\"/\"0\" -> w:val=\"true\"/\"false\"\n#     throughout styles.xml's latentStyles/lsdException block), plus one\n#     extLst child's namespace prefix being renamed (thm15 -> m) in\n#     theme1.xml.\n# No paragraph, run, bookmark, style definition, theme color, or\n# header/footer text is added, removed, or retitled anywhere in the\n# diff. That matches the commit message (\"fix import and update golden\n# files\"): a test-fixture refresh after fixing an *importer*, not a\n# content edit. The Word COM object model has no surface for dictating\n# XML namespace prefixes, attribute order, or boolean-literal spelling \u2014\n# the host always re-serializes OOXML with its own writer \u2014 so there is\n# nothing for a content-level script to change.\n#\n# We still touch the one piece of content the diff's context lines pass\n# through (the document's \"_GoBack\" bookmark) to confirm it is present\n# and unchanged, without performing any write that would introduce a new\n# difference (re-adding/moving the bookmark would duplicate or shift it,\n# which is not what this commit did).\n$d = $word.ActiveDocument\n\n$bookmarkExists = $d.Bookmarks.Exists(\"_GoBack\")\nWrite-Output \"GoBack bookmark present: $bookmarkExists\"\n\n# Read-only verification only -- the bookmark's name/location and every\n# paragraph/style/theme value already match the target state, so no\n# property is set and no content is mutated.\n"}
